$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    2 = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 5.582307763322248)
    3 = @(1.445647641019636, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 4.327115817150455)
    4 = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 5.582307763322248)
    5 = @(0.04172184405617529, 0.04103571897497393, 3.223369029078222, 13.86384647080068, 17.16997306291006)
    6 = @(0.2881169905109251, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 5.671859677732939)
}

foreach ($row in $values.Keys) {
    $rowValues = $values[$row]
    $ws.Range("B$row").Value = $rowValues[0]
    $ws.Range("C$row").Value = $rowValues[1]
    $ws.Range("D$row").Value = $rowValues[2]
    $ws.Range("E$row").Value = $rowValues[3]
    $ws.Range("G$row").Value = $rowValues[4]
}
